# repull data, push all data, mean calculation
# Updates the "dSF" column (F) values that changed after re-pulling the
# upstream data. Column E ("dS0") is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    10 = -3
    11 = -5
    12 = 8
    13 = 8
    18 = 4
    22 = -8
    24 = -5
    25 = 1
    26 = 1
    27 = -4
    28 = -10
    31 = -1
    35 = -2
    36 = -1
    37 = -2
    38 = -2
    41 = 3
    44 = -1
    45 = 0
    49 = 0
    52 = -2
    53 = -1
    54 = -1
    59 = 2
    62 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
